$d = $word.ActiveDocument

# --- 1. Trim the trailing "ABCDEFG20.75mm4mwforce-none16x1616x16" run group
#        (2nd & 3rd <w:br/><w:t> pairs) from the 3rd paragraph, keeping just
#        the leading <w:br/> + 24 spaces.
$p = $d.Paragraphs(3)
$pStart = $p.Range.Start
$fullText = $p.Range.Text
$vt = [char]11
$idx1 = $fullText.IndexOf($vt)
$idx2 = $fullText.IndexOf($vt, $idx1 + 1)
$delStart = $pStart + $idx2
$delEnd = $p.Range.End - 1
$r = $d.Range($delStart, $delEnd)
$r.Text = ""

# --- 2. Split the "Welcome to the Canara..." paragraph: the leading
#        "<w:br/>spaces" stays in its own paragraph, "Welcome to..." moves
#        to a new paragraph right after it.
$p = $d.Paragraphs(6)
$pStart = $p.Range.Start
$fullText = $p.Range.Text
$idx = $fullText.IndexOf("Welcome")
$insertPos = $pStart + $idx
$r = $d.Range($insertPos, $insertPos)
$r.InsertParagraphBefore()

# --- 3. Truncate the "We are confident..." paragraph so it ends right
#        after "within 1".
$ok = $d.Content.Find.Execute("within 15 days from receipt of this document.", $true, $false, $false, $false, $false, $true, 1, $false, "within 1", 2)

# --- 4. Insert two new leading paragraphs: "210mm" and "158000570".
$r = $d.Paragraphs(1).Range
$r.InsertParagraphBefore()
$d.Paragraphs(1).Range.Text = "210mm"
$r2 = $d.Paragraphs(2).Range
$r2.InsertParagraphBefore()
$d.Paragraphs(2).Range.Text = "158000570"

# --- 5. Remove the now-unused custom paragraph style "Body_Colour_Inline".
$s = $d.Styles("Body_Colour_Inline")
$s.Delete()
